# Add a new "ano" (year) column in column A, to the left of the existing
# "mes" (B) / "drogas_kg" (C) columns.
#
#   A2 = "ano"  (new header)
#   A3:A14 = 2023  (year value for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("A2").Value = "ano"

# Fill the year value for every existing data row (rows 3-14, matching
# jan..dez in column B).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = 2023
}
